$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.683.55"
$ws.Range("E2").Value = "  +7.56%  "
$ws.Range("D3").Value = "3.468.66"
$ws.Range("E3").Value = "  +5.05%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "415.27"
$ws.Range("E5").Value = "  +3.64%  "
$ws.Range("D6").Value = "125.73"
$ws.Range("E6").Value = "  +14.31%  "
$ws.Range("D7").Value = "3.462.31"
$ws.Range("E7").Value = "  +5.13%  "
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "0.681"
$ws.Range("E10").Value = "  +8.63%  "
$ws.Range("D11").Value = "0.129"
$ws.Range("E11").Value = "  +32.55%  "
$ws.Range("D12").Value = "41.56"
$ws.Range("E12").Value = "  +4.86%  "
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "4.014.57"
$ws.Range("E14").Value = "  +5.67%  "
$ws.Range("D15").Value = "8.63"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "20.05"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("D17").Value = "3.461.31"
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("D18").Value = "62.616.75"
$ws.Range("E18").Value = "  +8.30%  "
$ws.Range("D19").Value = "1.04"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "10.88"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "0.0000140"
$ws.Range("E21").Value = "  +28.92%  "
$ws.Range("D22").Value = "3.35"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "317.04"
$ws.Range("E23").Value = "  +5.25%  "
$ws.Range("D24").Value = "82.30"
$ws.Range("E24").Value = "  +9.98%  "
$ws.Range("D25").Value = "13.13"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "3.19"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "30.94"
$ws.Range("E27").Value = "  +8.93%  "
$ws.Range("D28").Value = "7.85"
$ws.Range("E28").Value = "  +5.84%  "
$ws.Range("D29").Value = "7.93"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "0.176"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  +3.66%  "
$ws.Range("E33").Value = "  +23.37%  "
$ws.Range("D34").Value = "11.61"
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("D35").Value = "42.51"
$ws.Range("E35").Value = "  +3.68%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").Value = "0.0498"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "52.31"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "3.52"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "3.05"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("D42").Value = "2.01"
$ws.Range("E42").Value = "  +5.98%  "
$ws.Range("E43").Value = "  +2.70%  "
$ws.Range("D44").Value = "136.34"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "0.287"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "16.99"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "3.93"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "2.24"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "22.03"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").Value = "2.207.90"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").Value = "2.45"
$ws.Range("E51").Value = "  -0.54%  "
